$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue $ws "D2" "26.209.58"
Set-TextValue $ws "E2" "  -1.91%  "
Set-TextValue $ws "D3" "1.582.53"
Set-TextValue $ws "E4" "  -0.44%  "
Set-TextValue $ws "D5" "209.60"
Set-TextValue $ws "E5" "  -1.00%  "
Set-TextValue $ws "E6" "  -3.00%  "
Set-TextValue $ws "E7" "  -0.37%  "
Set-TextValue $ws "D8" "0.0611"
Set-TextValue $ws "E8" "  -1.34%  "
Set-TextValue $ws "E9" "  -0.52%  "
Set-TextValue $ws "D11" "0.0847"
Set-TextValue $ws "E11" "  +0.09%  "
Set-TextValue $ws "D12" "1.805.01"
Set-TextValue $ws "E12" "  -1.14%  "
Set-TextValue $ws "D13" "1.596.84"
Set-TextValue $ws "E13" "  -0.17%  "
Set-TextValue $ws "D14" "4.05"
Set-TextValue $ws "E14" "  +0.27%  "
Set-TextValue $ws "E15" "  -1.40%  "
Set-TextValue $ws "D16" "64.48"
Set-TextValue $ws "E16" "  -0.71%  "
Set-TextValue $ws "D17" "26.200.98"
Set-TextValue $ws "E17" "  -1.83%  "
Set-TextValue $ws "D18" "0.0₃0735"
Set-TextValue $ws "E18" "  -0.68%  "
Set-TextValue $ws "D19" "7.29"
Set-TextValue $ws "E19" "  +1.43%  "
Set-TextValue $ws "E20" "  -0.36%  "
Set-TextValue $ws "D21" "207.15"
Set-TextValue $ws "E21" "  -1.60%  "
Set-TextValue $ws "D22" "4.26"
Set-TextValue $ws "D23" "2.19"
Set-TextValue $ws "E23" "  -3.33%  "
Set-TextValue $ws "D24" "8.88"
Set-TextValue $ws "E24" "  -0.84%  "
Set-TextValue $ws "D25" "144.93"
Set-TextValue $ws "E25" "  +0.61%  "
Set-TextValue $ws "E26" "  -0.54%  "
Set-TextValue $ws "E27" "  -0.65%  "
Set-TextValue $ws "E28" "  -0.98%  "
Set-TextValue $ws "D29" "15.22"
Set-TextValue $ws "E29" "  -0.89%  "
Set-TextValue $ws "E30" "  -1.25%  "
Set-TextValue $ws "D31" "1.15"
Set-TextValue $ws "E31" "  -1.14%  "
Set-TextValue $ws "E32" "  -1.34%  "
Set-TextValue $ws "E33" "  -0.74%  "
Set-TextValue $ws "D34" "1.282.79"
Set-TextValue $ws "E34" "  -0.68%  "
Set-TextValue $ws "D35" "2.47"
Set-TextValue $ws "E35" "  -0.32%  "
Set-TextValue $ws "E36" "  +6.62%  "
Set-TextValue $ws "D37" "0.609"
Set-TextValue $ws "E37" "  +2.02%  "
Set-TextValue $ws "E38" "  -0.74%  "
Set-TextValue $ws "E39" "  -1.46%  "
Set-TextValue $ws "D40" "0.814"
Set-TextValue $ws "E40" "  -1.71%  "
Set-TextValue $ws "E41" "  +3.45%  "
Set-TextValue $ws "D42" "0.766"
Set-TextValue $ws "E42" "  -1.96%  "
Set-TextValue $ws "E43" "  -2.94%  "
Set-TextValue $ws "D44" "62.34"
Set-TextValue $ws "E44" "  -1.05%  "
Set-TextValue $ws "D45" "1.718.32"
Set-TextValue $ws "E45" "  -1.22%  "
Set-TextValue $ws "D46" "88.94"
Set-TextValue $ws "E46" "  -1.88%  "
Set-TextValue $ws "E47" "  -0.07%  "
Set-TextValue $ws "D49" "0.0506"
Set-TextValue $ws "B50" "USDD"
Set-TextValue $ws "C50" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws "D50" "1.00"
Set-TextValue $ws "E50" "  -0.08%  "
Set-TextValue $ws "B51" "BabyDogeCoin"
Set-TextValue $ws "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D51" "0.0₇0950"
Set-TextValue $ws "E51" "  -10.66%  "
